$d = $word.ActiveDocument

# --- 1. Remove the "_GoBack" bookmark from its current position (start of doc) ---
try {
    $gb = $d.Bookmarks.Item("_GoBack")
    $gb.Delete()
} catch {
}

# --- 2. Trim the trailing underscore "blank lines" off the label paragraphs ---
# Each old string is unique in the document, so a literal (non-wildcard) Find &
# Replace is unambiguous.  wdReplace numeric literals: wdReplaceOne = 1.
$edits = @(
    @("Cover Letter Author:  ____________________________________________________________", "Cover Letter Author:  "),
    @("Peer Reviewer:  _________________________________________________________________", "Peer Reviewer:  "),
    @("Meets expectations of genre?  _____________________________________________________", "Meets expectations of genre?  "),
    @("What is most important, interesting or striking?  ______________________________________", "What is most important, interesting or striking?  "),
    @("What is almost but not quite said?  _________________________________________________", "What is almost but not quite said?  "),
    @("Strengths based on the rubric:  ____________________________________________________", "Strengths based on the rubric:  "),
    @("Weaknesses based on the rubric:  __________________________________________________", "Weaknesses based on the rubric:  ")
)

foreach ($pair in $edits) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Range(0, $d.Content.End)
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# --- 3. Delete every paragraph whose entire content is the 80-underscore
# placeholder rule.  We locate each occurrence with Find (which gives us
# reliable Start/End offsets), then delete a fresh Range that also swallows
# the trailing paragraph mark so the paragraph itself collapses away
# (deleting a Range obtained straight off Find only clears the run text and
# leaves an empty paragraph behind, so we rebuild the Range from the
# Start/End numbers instead).
$underscoreRun = "______________________________________________________________________________"
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute($underscoreRun, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $s = $rng.Start
    $e = $rng.End
    $delRange = $d.Range($s, $e + 1)
    $delRange.Delete()
    $searchStart = $s
}

# --- 4. Insert the extra blank paragraphs that appear between sections.
# The sections all keep at least one blank paragraph already, so we just
# grow that run to the required length by inserting right after the label
# paragraph itself (all the blank paragraphs are identical/unformatted, so
# where exactly inside the run they land is not visible in the output).
$insertsAfterLabel = @(
    @("Meets expectations of genre?  ", 1),
    @("What is most important, interesting or striking?  ", 1),
    @("What is almost but not quite said?  ", 1),
    @("Strengths based on the rubric:  ", 0),
    @("Weaknesses based on the rubric:  ", 2)
)

foreach ($pair in $insertsAfterLabel) {
    $label = $pair[0]
    $count = $pair[1]
    if ($count -eq 0) { continue }
    $rng = $d.Range(0, $d.Content.End)
    $found = $rng.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $target = $d.Range($rng.Start, $rng.End)
        for ($i = 0; $i -lt $count; $i++) {
            $target.InsertParagraphAfter()
        }
    }
}

# --- 5. Re-anchor the "_GoBack" bookmark on the last (now empty) paragraph. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null

Write-Output "done"
